$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("X2:X20")
$rng.Font.ThemeColor = 1
$rng.Interior.ColorIndex = 3
$rng.Interior.ColorIndex = -4142
Write-Output "done"
